$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K") for rows 2 through 38,
# regenerated per commit: "use K instead of Strike#"
$newValues = @(3,4,6,5,8,6,8,2,2,6,3,7,6,5,3,4,5,7,4,3,5,4,3,0,4,3,4,3,5,2,6,2,2,6,2,4,4)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
